# Update res_bus/vm_pu.xlsx results for the 380 kV case (Case_4_161)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$vals = @(1.02, 1.041996964640007, 1.043860229030055, 1.055527091839851, 1.062840747374948, 1.039351185376353, 1.047075161650013, 1.046632905581689, 1.058267305117526, 1.065561003750701, 1.048562128480756)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "2").Value = $vals[$i] }

$vals = @(1.02, 1.043147733651574, 1.044734762056177, 1.056663638857233, 1.064045320553874, 1.039626729482457, 1.047871187048731, 1.047318393646616, 1.059216487178757, 1.066579492489757, 1.049359284326846)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "3").Value = $vals[$i] }

$vals = @(1.02, 1.043892177384048, 1.0453001034771, 1.05739921494096, 1.064824921065073, 1.039803158006219, 1.048385551134136, 1.047760774462567, 1.059830241114861, 1.067238111428524, 1.04987437886825)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "4").Value = $vals[$i] }

$vals = @(1.02, 1.044205099550369, 1.045537643954096, 1.057708488856159, 1.065152704660421, 1.03987688168705, 1.048601618595931, 1.047946470189445, 1.060088160964633, 1.067514897143637, 1.050090753170624)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "5").Value = $vals[$i] }

$vals = @(1.02, 1.044257638135039, 1.045577520483247, 1.057760419543711, 1.065207743315164, 1.03988923401494, 1.048637887251336, 1.047977632840453, 1.060131460863849, 1.067561364977856, 1.05012707333168)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "6").Value = $vals[$i] }

$vals = @(1.02, 1.04389635882868, 1.045303278011448, 1.05740334732488, 1.064829300768948, 1.039804144861463, 1.048388438910108, 1.047763256842375, 1.059833687853622, 1.067241810235652, 1.049877270745196)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "7").Value = $vals[$i] }

$vals = @(1.02, 1.042385910626383, 1.044155893480695, 1.055911161801319, 1.063247805744872, 1.039444693407699, 1.047344331095822, 1.04686481313495, 1.058588175174088, 1.065905292449777, 1.048831680178052)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "8").Value = $vals[$i] }

$vals = @(1.02, 1.039722864868481, 1.042129910146776, 1.053282883492173, 1.060462199372533, 1.038796991441749, 1.045498956938774, 1.045272620184609, 1.056390094185703, 1.063546990247883, 1.046983685378034)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "9").Value = $vals[$i] }

$vals = @(1.02, 1.03794644612953, 1.040776450158389, 1.05153141103122, 1.058605863676277, 1.03835556624557, 1.04426495387225, 1.044205069395647, 1.054922420301131, 1.061972591960235, 1.045747929885613)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "10").Value = $vals[$i] }

$vals = @(1.02, 1.03717696728285, 1.040189716919136, 1.050773159906031, 1.057802208780314, 1.03816213928267, 1.043729716785978, 1.043741357673291, 1.054286347230106, 1.061290326093245, 1.045211932701285)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "11").Value = $vals[$i] }

$vals = @(1.02, 1.036891105035091, 1.039971675941289, 1.050491532652413, 1.057503716767341, 1.038089948023021, 1.043530768894195, 1.043568895217366, 1.054049996028222, 1.061036819658991, 1.045012702280678)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "12").Value = $vals[$i] }

$vals = @(1.02, 1.036952425457081, 1.040018451101606, 1.050551941725463, 1.057567743381388, 1.038105448858916, 1.043573450107447, 1.043605898932797, 1.054100698043292, 1.061091201403435, 1.045055444106148)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "13").Value = $vals[$i] }

$vals = @(1.02, 1.037153338718999, 1.040171695670475, 1.050749880088818, 1.057777534913349, 1.038156178947212, 1.043713274491576, 1.043727106349824, 1.054266812113478, 1.061269372862701, 1.045195467056939)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "14").Value = $vals[$i] }

$vals = @(1.02, 1.037277122186673, 1.040266101157527, 1.050871839188878, 1.057906797165817, 1.038187389853717, 1.043799406726272, 1.043801757164097, 1.054369149132324, 1.061379139208384, 1.045281721609289)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "15").Value = $vals[$i] }

$vals = @(1.02, 1.037997507899268, 1.04081537541928, 1.051581736672032, 1.058659202657116, 1.038368355148528, 1.044300456633847, 1.044235813701037, 1.054964622469542, 1.062017860195241, 1.045783483065204)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "16").Value = $vals[$i] }

$vals = @(1.02, 1.038449311174764, 1.041159739249056, 1.052027075346342, 1.059131206076088, 1.03848125741772, 1.044614508877417, 1.044507695848162, 1.055337996287243, 1.062418367282738, 1.046097981298982)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "17").Value = $vals[$i] }

$vals = @(1.02, 1.03871281369823, 1.041360535645009, 1.052286848226192, 1.059406532323388, 1.038546890764126, 1.044797602922942, 1.044666139744399, 1.055555724904395, 1.062651924173334, 1.046281335359052)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "18").Value = $vals[$i] }

$vals = @(1.02, 1.038802656778101, 1.041428990984473, 1.052375426540008, 1.059500413985449, 1.038569232623898, 1.044860018479531, 1.044720141230581, 1.055629955608514, 1.062731552234697, 1.046343839552891)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "19").Value = $vals[$i] }

$vals = @(1.02, 1.038400839731863, 1.041122799016828, 1.051979293244191, 1.059080563072252, 1.03846916689622, 1.044580823081874, 1.044478539989482, 1.055297942404934, 1.062375402062024, 1.046064247665746)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "20").Value = $vals[$i] }

$vals = @(1.02, 1.037094175987523, 1.040126571788149, 1.050691591605432, 1.057715755962566, 1.038141249695533, 1.043672103485175, 1.043691419844861, 1.054217898024974, 1.06121690809645, 1.045154237582987)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "21").Value = $vals[$i] }

$vals = @(1.02, 1.036272371805629, 1.039499613027938, 1.049882083204363, 1.056857769489249, 1.03793308544633, 1.043099962256642, 1.043195257181561, 1.053538336934318, 1.060488038948345, 1.044581283848285)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "22").Value = $vals[$i] }

$vals = @(1.02, 1.036708050320828, 1.039832031949668, 1.050311207903649, 1.057312593279627, 1.038043625950058, 1.043403340647012, 1.043458402901692, 1.053898632357813, 1.060874471997528, 1.044885093070767)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "23").Value = $vals[$i] }

$vals = @(1.02, 1.038422741972414, 1.041139490918891, 1.052000883874601, 1.059103446420099, 1.038474630764036, 1.044596044512355, 1.04449171470044, 1.055316041196772, 1.062394816356478, 1.046079490712403)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "24").Value = $vals[$i] }

$vals = @(1.02, 1.040411505737277, 1.042654168387399, 1.05396222665196, 1.061182211026842, 1.038966133620233, 1.045976689008801, 1.045685311630374, 1.056958750240212, 1.064157051116754, 1.047462095882397)
for ($i = 0; $i -lt $cols.Length; $i++) { $ws.Range($cols[$i] + "25").Value = $vals[$i] }

Write-Output "Updated 380 kV case values (rows 2-25)"